$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after the last data row (row 20), pushing the
# blank gap + footer rows (signature block) down from 25/26 to 27/28.
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

# Duplicate the last two worker rows (19:20 -> 21:22) including formatting,
# mirroring "EC" rows being carried forward for a new period.
$ws.Range("B19:J20").Copy()
$ws.Range("B21").PasteSpecial()

# New rows belong to period 2509 (the newly added period).
$ws.Range("E21").Value = "2509"
$ws.Range("E22").Value = "2509"

# Update summary figures: one more period counted, and VALOR MORA grows by
# the Valor Mora of the two newly added rows (56940 + 56940).
$ws.Range("F13").Value = 5
$ws.Range("E11").Value = 327135
